$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 2 de Agosto de 2020 a las 21:03"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4791811
$ws.Range("C4").Value = 27493
$ws.Range("D4").Value = 2372808
$ws.Range("E4").Value = 2260815
$ws.Range("G4").Value = 290
$ws.Range("H4").Value = 158188

# Row 6 - India
$ws.Range("B6").Value = 1804702
$ws.Range("C6").Value = 52783
$ws.Range("D6").Value = 1187228
$ws.Range("E6").Value = 579313
$ws.Range("G6").Value = 758
$ws.Range("H6").Value = 38161

# Row 21 - Alemania
$ws.Range("B21").Value = 211436
$ws.Range("C21").Value = 359
$ws.Range("E21").Value = 8610

# Row 56 - Afganistan
$ws.Range("D56").Value = 25510
$ws.Range("E56").Value = 9916
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 1284

# Row 118 - Sri Lanka
$ws.Range("B118").Value = 2823
$ws.Range("C118").Value = 8
$ws.Range("E118").Value = 298

# Row 123 - Sudan del Sur
$ws.Range("B123").Value = 2429
$ws.Range("C123").Value = 77
$ws.Range("E123").Value = 1208

# Row 125 - Namibia
$ws.Range("B125").Value = 2294
$ws.Range("C125").Value = 70
$ws.Range("D125").Value = 187
$ws.Range("E125").Value = 2096

# Rows 155-156: swap Botsuana/Siria order and update their data
$ws.Range("A155").Value = "Siria"
$ws.Range("B155").Value = 809
$ws.Range("C155").Value = 29
$ws.Range("D155").Value = 256
$ws.Range("E155").Value = 509
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 44

$ws.Range("A156").Value = "Botsuana"
$ws.Range("B156").Value = 804
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 63
$ws.Range("E156").Value = 739
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 2
